# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# on Sheet1 of the cryptos workbook to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.562.29"
$ws.Range("E2").Value = "  +2.06%  "
$ws.Range("D3").Value = "2.291.41"
$ws.Range("E3").Value = "  +1.12%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.11"
$ws.Range("E5").Value = "  +1.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.04"
$ws.Range("E6").Value = "  +6.46%  "
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.498"
$ws.Range("E9").Value = "  +3.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.39"
$ws.Range("E10").Value = "  +12.68%  "
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("E12").Value = "  -1.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.75"
$ws.Range("E13").Value = "  +2.46%  "
$ws.Range("D14").Value = "2.644.37"
$ws.Range("E14").Value = "  +0.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.63"
$ws.Range("E15").Value = "  +3.11%  "
$ws.Range("D16").Value = "2.299.35"
$ws.Range("E16").Value = "  +1.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.810"
$ws.Range("E17").Value = "  +6.05%  "
$ws.Range("D18").Value = "42.440.64"
$ws.Range("E18").Value = "  +1.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.81"
$ws.Range("E19").Value = "  +1.42%  "
$ws.Range("E20").Value = "  +1.94%  "
$ws.Range("E21").Value = "  +2.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.01"
$ws.Range("E22").Value = "  +1.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "243.16"
$ws.Range("E23").Value = "  +1.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.61"
$ws.Range("E24").Value = "  +1.07%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.04"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.76"
$ws.Range("E28").Value = "  +9.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.62"
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("E30").Value = "  +2.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "161.08"
$ws.Range("E31").Value = "  +0.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.32"
$ws.Range("E32").Value = "  +1.76%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("E34").Value = "  +4.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0754"
$ws.Range("E35").Value = "  +1.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.40"
$ws.Range("E36").Value = "  +2.96%  "
$ws.Range("E37").Value = "  +3.75%  "
$ws.Range("E38").Value = "  +5.00%  "
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("E40").Value = "  -0.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.22"
$ws.Range("E41").Value = "  +6.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.42"
$ws.Range("E42").Value = "  +16.94%  "
$ws.Range("D43").Value = "2.005.90"
$ws.Range("E43").Value = "  -1.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.36"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("E45").Value = "  +3.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.04"
$ws.Range("E46").Value = "  +5.58%  "
$ws.Range("E47").Value = "  -1.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.89"
$ws.Range("E48").Value = "  +3.97%  "
$ws.Range("E49").Value = "  +1.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.03"
$ws.Range("E50").Value = "  +0.49%  "
$ws.Range("E51").Value = "  -0.60%  "
